# Apply crypto price/volume updates per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.703.59'
$ws.Range("E2").Value = '  -1.10%  '

$ws.Range("D3").Value = '3.076.42'
$ws.Range("E3").Value = '  -2.63%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.46%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.74'
$ws.Range("E5").Value = '  -0.60%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.86'
$ws.Range("E6").Value = '  +3.42%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").Value = '3.071.72'
$ws.Range("E9").Value = '  -2.77%  '

$ws.Range("E10").Value = '  -2.81%  '

$ws.Range("E11").Value = '  -0.07%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.456'
$ws.Range("E12").Value = '  -1.63%  '

$ws.Range("E13").Value = '  -2.74%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.53'
$ws.Range("E14").Value = '  +0.61%  '

$ws.Range("D15").Value = '3.586.40'
$ws.Range("E15").Value = '  -2.80%  '

$ws.Range("E16").Value = '  -2.30%  '

$ws.Range("B17").Value = 'Polkadot'
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.15'
$ws.Range("E17").Value = '  -0.99%  '

$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '63.574.62'
$ws.Range("E18").Value = '  -0.93%  '

$ws.Range("D19").Value = '3.071.75'
$ws.Range("E19").Value = '  -2.66%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '469.01'
$ws.Range("E20").Value = '  +0.06%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.62'
$ws.Range("E21").Value = '  +1.17%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.726'
$ws.Range("E22").Value = '  -1.70%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.53'
$ws.Range("E23").Value = '  +0.73%  '

$ws.Range("B24").Value = 'InternetComputer(DFINITY)'
$ws.Range("C24").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.24'
$ws.Range("E24").Value = '  +1.05%  '

$ws.Range("B25").Value = 'Fetch.AI'
$ws.Range("C25").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.37'
$ws.Range("E25").Value = '  +0.32%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '81.10'
$ws.Range("E26").Value = '  -0.23%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.09%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.85'
$ws.Range("E28").Value = '  +0.59%  '

$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.68'
$ws.Range("E29").Value = '  -1.56%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.42%  '

$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.30'
$ws.Range("E31").Value = '  +0.18%  '

$ws.Range("E32").Value = '  -1.55%  '

$ws.Range("E33").Value = '  +4.77%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '27.34'
$ws.Range("E34").Value = '  -0.93%  '

$ws.Range("D35").Value = '0.0₃0849'
$ws.Range("E35").Value = '  -0.78%  '

$ws.Range("E36").Value = '  -1.10%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.37'
$ws.Range("E37").Value = '  +2.23%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.11'
$ws.Range("E38").Value = '  -0.57%  '

$ws.Range("E39").Value = '  -5.08%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '9.29'
$ws.Range("E40").Value = '  +1.57%  '

$ws.Range("E41").Value = '  -2.87%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '447.52'
$ws.Range("E42").Value = '  -2.05%  '

$ws.Range("E43").Value = '  -2.44%  '

$ws.Range("E44").Value = '  -2.37%  '

$ws.Range("D45").Value = '2.830.58'
$ws.Range("E45").Value = '  -3.32%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.89'
$ws.Range("E46").Value = '  -1.20%  '

$ws.Range("E47").Value = '  -0.30%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '129.27'
$ws.Range("E48").Value = '  +1.48%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.46'
$ws.Range("E49").Value = '  +3.49%  '

$ws.Range("E50").Value = '  +0.00%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.25'
$ws.Range("E51").Value = '  -0.59%  '
